$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nationalité")

# Remove the "436 / BAHAMAS" nationality entry (row 168) by shifting all the
# rows below it up by one position, then deleting the now-duplicate last data
# row so the trailing blank spacer row becomes row 183.
for ($r = 168; $r -le 182; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r + 1, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r + 1, 3).Value2
}

$ws.Rows.Item(183).Delete()
